# Update justification (column C) labels: "Postulado 1/2", "Proposición I.#",
# "Axioma 5" are renamed/renumbered with Roman numerals, and the equality
# statements in column B are switched from "=" to the congruence symbol
# "\cong". Step 23's narrative text is also reworded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (justificacion) relabeling ---
$ws.Range("C8").Value  = "Postulado II"
$ws.Range("C10").Value = "Proposición I.X"
$ws.Range("C11").Value = "Postulado I"
$ws.Range("C12").Value = "Postulado I"
$ws.Range("C13").Value = "Postulado I"
$ws.Range("C14").Value = "Postulado II"
$ws.Range("C16").Value = "Postulado I"
$ws.Range("C17").Value = "Proposición I.III"
$ws.Range("C19").Value = "Postulado I"
$ws.Range("C20").Value = "Postulado I"
$ws.Range("C21").Value = "Proposición I.X"
$ws.Range("C22").Value = "Proposición I.III"
$ws.Range("C23").Value = "Proposición I.XV"
$ws.Range("C24").Value = "Proposición I.XV"
$ws.Range("C25").Value = "Proposición I.IV"
$ws.Range("C26").Value = "Axioma X"
$ws.Range("C27").Value = "Axioma X"
$ws.Range("C28").Value = "Proposición I.IV"

# --- Column B (descripcion) wording changes: "=" -> "\cong" (congruence) ---
$ws.Range("B21").Value = '$\overline{AE} \cong \overline{EC}$ '
$ws.Range("B22").Value = '$\overline{BE} \cong \overline{EG}$ '
$ws.Range("B23").Value = '$\angle BEA \cong \angle GEC$'
$ws.Range("B24").Value = '$\angle BEA \cong \angle GEC$'
$ws.Range("B25").Value = 'De las afirmaciones (20), (22), (21), se deduce que $\triangle ABE \cong \triangle FEC$'
$ws.Range("B28").Value = '$\angle BAE \cong \angle ECG$ '

# --- Cell E10 loses its custom (applyFont) style, reverting to the default ---
$ws.Range("E10").Style = "Normal"

# --- Sheet view: scrolled position / active selection changed ---
$ws.Range("B22").Select()
